# Adds the new "7.0 Login with locked user" test case (rows 27-30) to the
# "test_Login_ValidData" sheet, mirroring the structure/formatting of the
# immediately preceding "6.0" test case (rows 23-26).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test_Login_ValidData")

function Copy-CellFormat {
    param($srcAddr, $dstAddr)
    $src = $ws.Range($srcAddr)
    $dst = $ws.Range($dstAddr)
    $dst.NumberFormat = $src.NumberFormat
    $dst.Font.Name = $src.Font.Name
    $dst.Font.Size = $src.Font.Size
    $dst.Font.Bold = $src.Font.Bold
    $dst.Font.Underline = $src.Font.Underline
    $dst.Font.Color = $src.Font.Color
    $dst.Interior.Color = $src.Interior.Color
    $dst.HorizontalAlignment = $src.HorizontalAlignment
    $dst.VerticalAlignment = $src.VerticalAlignment
    $dst.WrapText = $src.WrapText
}

# Row 27 mirrors row 23 (section header row: step number + title banner)
Copy-CellFormat "A23" "A27"
Copy-CellFormat "B23" "B27"
Copy-CellFormat "C23" "C27"
Copy-CellFormat "D23" "D27"
Copy-CellFormat "E23" "E27"
Copy-CellFormat "F23" "F27"

# Row 28 mirrors row 24
Copy-CellFormat "A24" "A28"
Copy-CellFormat "B24" "B28"
Copy-CellFormat "C24" "C28"
Copy-CellFormat "D24" "D28"
Copy-CellFormat "E24" "E28"
Copy-CellFormat "F24" "F28"

# Row 29 mirrors row 25
Copy-CellFormat "A25" "A29"
Copy-CellFormat "B25" "B29"
Copy-CellFormat "C25" "C29"
Copy-CellFormat "D25" "D29"
Copy-CellFormat "E25" "E29"
Copy-CellFormat "F25" "F29"

# Row 30 mirrors row 26
Copy-CellFormat "A26" "A30"
Copy-CellFormat "B26" "B30"
Copy-CellFormat "C26" "C30"
Copy-CellFormat "D26" "D30"
Copy-CellFormat "E26" "E30"
Copy-CellFormat "F26" "F30"

# ---- Values ----

# Row 27: new test case header "7.0 - Login with locked user"
$ws.Range("A27").Value2 = 7
$ws.Range("B27").Value2 = "Login with locked user"

# Row 28: step 6.1
$ws.Range("A28").NumberFormat = "@"
$ws.Range("A28").Value2 = "6.1"
$ws.Range("B28").Value2 = "Enter a valid username"
$ws.Range("C28").Value2 = "The username has been inserted"
$ws.Range("D28").Value2 = "locked_out_user"
$ws.Range("E28").Value2 = "The username field is empty"

# Row 29: step 6.2
$ws.Range("A29").NumberFormat = "@"
$ws.Range("A29").Value2 = "6.2"
$ws.Range("B29").Value2 = "Enter a valid passowrd"
$ws.Range("C29").Value2 = "The password has been inserted"
$ws.Range("D29").Value2 = "secret_sauce"
$ws.Range("E29").Value2 = "The password field is empty"

# Row 30: step 6.3
$ws.Range("A30").NumberFormat = "@"
$ws.Range("A30").Value2 = "6.3"
$ws.Range("B30").Value2 = "Click on Login button"
$ws.Range("C30").Value2 = 'Login is unsuccessfull, user is not being redirected to home page. Error "Epic sadface: Sorry, this user has been locked out ".'
$ws.Range("E30").Value2 = 'Error "Epic sadface: Sorry, this user has been locked out ".'
